$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A165").NumberFormat = "@"
$ws.Range("A165").Value = "2026-02-14 12:36:42"
$ws.Range("A165").Style = "Normal"
$ws.Range("B165").NumberFormat = "@"
$ws.Range("B165").Value = "237681102046"
$ws.Range("B165").Style = "Normal"
$ws.Range("C165").Value = "FRANCOISE NKENFACK NKENGMO"
$ws.Range("D165").Value = 397988
$ws.Range("A166").NumberFormat = "@"
$ws.Range("A166").Value = "2026-02-15 12:54:57"
$ws.Range("A166").Style = "Normal"
$ws.Range("B166").NumberFormat = "@"
$ws.Range("B166").Value = "237682368679"
$ws.Range("B166").Style = "Normal"
$ws.Range("C166").Value = "MFS SIM AA 2"
$ws.Range("D166").Value = 1493
$ws.Range("A167").NumberFormat = "@"
$ws.Range("A167").Value = "2026-02-15 18:13:33"
$ws.Range("A167").Style = "Normal"
$ws.Range("B167").NumberFormat = "@"
$ws.Range("B167").Value = "237683360459"
$ws.Range("B167").Style = "Normal"
$ws.Range("C167").Value = "LUCIE MAJOLIE LELE NKANKEU"
$ws.Range("D167").Value = 1950
$ws.Range("A168").NumberFormat = "@"
$ws.Range("A168").Value = "2026-02-15 12:30:55"
$ws.Range("A168").Style = "Normal"
$ws.Range("B168").NumberFormat = "@"
$ws.Range("B168").Value = "237652899422"
$ws.Range("B168").Style = "Normal"
$ws.Range("C168").Value = "NOUMOU epouse SAGNON MARCELINE LA NEGRESSE"
$ws.Range("D168").Value = 549111
$ws.Range("A169").NumberFormat = "@"
$ws.Range("A169").Value = "2026-02-15 14:29:08"
$ws.Range("A169").Style = "Normal"
$ws.Range("B169").NumberFormat = "@"
$ws.Range("B169").Value = "237670904526"
$ws.Range("B169").Style = "Normal"
$ws.Range("C169").Value = "MFS SIM PROVISOIRE 20"
$ws.Range("D169").Value = 0
$ws.Range("A170").NumberFormat = "@"
$ws.Range("A170").Value = "2026-02-15 10:21:19"
$ws.Range("A170").Style = "Normal"
$ws.Range("B170").NumberFormat = "@"
$ws.Range("B170").Value = "237671105116"
$ws.Range("B170").Style = "Normal"
$ws.Range("C170").Value = "MFS  AM FACE HÔPITAL GÉNÉRAL"
$ws.Range("D170").Value = 30
$ws.Range("A171").NumberFormat = "@"
$ws.Range("A171").Value = "2026-02-15 18:06:17"
$ws.Range("A171").Style = "Normal"
$ws.Range("B171").NumberFormat = "@"
$ws.Range("B171").Value = "237672916354"
$ws.Range("B171").Style = "Normal"
$ws.Range("C171").Value = "MAMADOU DIAN BAH LENA GLOBAL"
$ws.Range("D171").Value = 32250
$ws.Range("A172").NumberFormat = "@"
$ws.Range("A172").Value = "2026-02-15 15:55:44"
$ws.Range("A172").Style = "Normal"
$ws.Range("B172").NumberFormat = "@"
$ws.Range("B172").Value = "237672920086"
$ws.Range("B172").Style = "Normal"
$ws.Range("C172").Value = "NAMY NGOKO CLARISSE ROSE VERTINE KAMILAH CONNECTION"
$ws.Range("D172").Value = 11730
$ws.Range("A173").NumberFormat = "@"
$ws.Range("A173").Value = "2026-02-15 08:10:25"
$ws.Range("A173").Style = "Normal"
$ws.Range("B173").NumberFormat = "@"
$ws.Range("B173").Value = "237674000053"
$ws.Range("B173").Style = "Normal"
$ws.Range("C173").Value = "YVETTE LAURE NGANDO SIMO EPSE DIMO"
$ws.Range("D173").Value = 237
$ws.Range("A174").NumberFormat = "@"
$ws.Range("A174").Value = "2026-02-15 14:32:51"
$ws.Range("A174").Style = "Normal"
$ws.Range("B174").NumberFormat = "@"
$ws.Range("B174").Value = "237674841555"
$ws.Range("B174").Style = "Normal"
$ws.Range("C174").Value = "BEATRICE TCHAMTIEU EPSE NGAMENI"
$ws.Range("D174").Value = 90817
$ws.Range("A175").NumberFormat = "@"
$ws.Range("A175").Value = "2026-02-15 14:43:18"
$ws.Range("A175").Style = "Normal"
$ws.Range("B175").NumberFormat = "@"
$ws.Range("B175").Value = "237674899678"
$ws.Range("B175").Style = "Normal"
$ws.Range("C175").Value = "VIVIANE MADJUIMEKEM FOMEKONG"
$ws.Range("D175").Value = 175766
$ws.Range("A176").NumberFormat = "@"
$ws.Range("A176").Value = "2026-02-15 10:20:03"
$ws.Range("A176").Style = "Normal"
$ws.Range("B176").NumberFormat = "@"
$ws.Range("B176").Value = "237676439452"
$ws.Range("B176").Style = "Normal"
$ws.Range("C176").Value = "RONIS BRAVO DONGMO TSAGUE"
$ws.Range("D176").Value = 219279
$ws.Range("A177").NumberFormat = "@"
$ws.Range("A177").Value = "2026-02-15 17:35:45"
$ws.Range("A177").Style = "Normal"
$ws.Range("B177").NumberFormat = "@"
$ws.Range("B177").Value = "237676695935"
$ws.Range("B177").Style = "Normal"
$ws.Range("C177").Value = "Mathieu Djitouo"
$ws.Range("D177").Value = 51303
$ws.Range("A178").NumberFormat = "@"
$ws.Range("A178").Value = "2026-02-15 14:56:28"
$ws.Range("A178").Style = "Normal"
$ws.Range("B178").NumberFormat = "@"
$ws.Range("B178").Value = "237677745809"
$ws.Range("B178").Style = "Normal"
$ws.Range("C178").Value = "JOSEPHINE CLAIRE NGUENKAM KADJI"
$ws.Range("D178").Value = 150103
$ws.Range("A179").NumberFormat = "@"
$ws.Range("A179").Value = "2026-02-15 16:00:03"
$ws.Range("A179").Style = "Normal"
$ws.Range("B179").NumberFormat = "@"
$ws.Range("B179").Value = "237679127464"
$ws.Range("B179").Style = "Normal"
$ws.Range("C179").Value = "DJUFFO TSOATA MARIE NOEL KAMILAH CONNECTION GROUP"
$ws.Range("D179").Value = 26352
$ws.Range("A180").NumberFormat = "@"
$ws.Range("A180").Value = "2026-02-15 14:31:06"
$ws.Range("A180").Style = "Normal"
$ws.Range("B180").NumberFormat = "@"
$ws.Range("B180").Value = "237679422291"
$ws.Range("B180").Style = "Normal"
$ws.Range("C180").Value = "ETS LE CONTENT 32"
$ws.Range("D180").Value = 100000
$ws.Range("A181").NumberFormat = "@"
$ws.Range("A181").Value = "2026-02-15 05:47:57"
$ws.Range("A181").Style = "Normal"
$ws.Range("B181").NumberFormat = "@"
$ws.Range("B181").Value = "237651433330"
$ws.Range("B181").Style = "Normal"
$ws.Range("C181").Value = "NGUIAZONG DORIANE LAURE KAMILAH CONNECTION GROUP"
$ws.Range("D181").Value = 59822
$ws.Range("A182").NumberFormat = "@"
$ws.Range("A182").Value = "2026-02-15 10:11:16"
$ws.Range("A182").Style = "Normal"
$ws.Range("B182").NumberFormat = "@"
$ws.Range("B182").Value = "237654168696"
$ws.Range("B182").Style = "Normal"
$ws.Range("C182").Value = "DZEUMAZONG FLORENCE ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Range("D182").Value = 1743
$ws.Range("A183").NumberFormat = "@"
$ws.Range("A183").Value = "2026-02-15 03:21:15"
$ws.Range("A183").Style = "Normal"
$ws.Range("B183").NumberFormat = "@"
$ws.Range("B183").Value = "237670799877"
$ws.Range("B183").Style = "Normal"
$ws.Range("C183").Value = "ETS MOBILE FINANCIAL SERVICES MFS MENANDJIO HORTENSE BIENVENUE"
$ws.Range("D183").Value = 681193
$ws.Range("A184").NumberFormat = "@"
$ws.Range("A184").Value = "2026-02-15 14:42:37"
$ws.Range("A184").Style = "Normal"
$ws.Range("B184").NumberFormat = "@"
$ws.Range("B184").Value = "237671351291"
$ws.Range("B184").Style = "Normal"
$ws.Range("C184").Value = "MFS LTDLA CBOX R3 MOUTHIEU JOSETTE CHANCELINE"
$ws.Range("D184").Value = 126983
$ws.Range("A185").NumberFormat = "@"
$ws.Range("A185").Value = "2026-02-15 15:32:08"
$ws.Range("A185").Style = "Normal"
$ws.Range("B185").NumberFormat = "@"
$ws.Range("B185").Value = "237671378136"
$ws.Range("B185").Style = "Normal"
$ws.Range("C185").Value = "KOUBINOM DIPITA SARIETTE CRISTELLE ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Range("D185").Value = 210112
$ws.Range("A186").NumberFormat = "@"
$ws.Range("A186").Value = "2026-02-15 15:08:40"
$ws.Range("A186").Style = "Normal"
$ws.Range("B186").NumberFormat = "@"
$ws.Range("B186").Value = "237671605749"
$ws.Range("B186").Style = "Normal"
$ws.Range("C186").Value = "ETS TCHATCHOUANG PAUL  ETP LTDLA CBOX RO MEGAPTCHE VICTORINE"
$ws.Range("D186").Value = 542405
$ws.Range("A187").NumberFormat = "@"
$ws.Range("A187").Value = "2026-02-15 11:32:13"
$ws.Range("A187").Style = "Normal"
$ws.Range("B187").NumberFormat = "@"
$ws.Range("B187").Value = "237671615641"
$ws.Range("B187").Style = "Normal"
$ws.Range("C187").Value = "BEGO FOGUE CHRISTELLE KAMILAH CONNECTION GROUP"
$ws.Range("D187").Value = 116
$ws.Range("A188").NumberFormat = "@"
$ws.Range("A188").Value = "2026-02-15 13:55:51"
$ws.Range("A188").Style = "Normal"
$ws.Range("B188").NumberFormat = "@"
$ws.Range("B188").Value = "237673739931"
$ws.Range("B188").Style = "Normal"
$ws.Range("C188").Value = "MOFFO GERMAIN SPECTRUM SPECTRUM"
$ws.Range("D188").Value = 105981
$ws.Range("A189").NumberFormat = "@"
$ws.Range("A189").Value = "2026-02-15 13:45:41"
$ws.Range("A189").Style = "Normal"
$ws.Range("B189").NumberFormat = "@"
$ws.Range("B189").Value = "237674466307"
$ws.Range("B189").Style = "Normal"
$ws.Range("C189").Value = "MFS MEYIMDOUNG MARIE LOUISE"
$ws.Range("D189").Value = 561158
$ws.Range("A190").NumberFormat = "@"
$ws.Range("A190").Value = "2026-02-15 14:36:03"
$ws.Range("A190").Style = "Normal"
$ws.Range("B190").NumberFormat = "@"
$ws.Range("B190").Value = "237674959564"
$ws.Range("B190").Style = "Normal"
$ws.Range("C190").Value = "FOTSO PASCAL KAMILAH CONNECTION GROUP"
$ws.Range("D190").Value = 1142
$ws.Range("A191").NumberFormat = "@"
$ws.Range("A191").Value = "2026-02-15 14:21:14"
$ws.Range("A191").Style = "Normal"
$ws.Range("B191").NumberFormat = "@"
$ws.Range("B191").Value = "237675208505"
$ws.Range("B191").Style = "Normal"
$ws.Range("C191").Value = "N A SAIDOU INOUSSA"
$ws.Range("D191").Value = 53324
$ws.Range("A192").NumberFormat = "@"
$ws.Range("A192").Value = "2026-02-15 16:11:48"
$ws.Range("A192").Style = "Normal"
$ws.Range("B192").NumberFormat = "@"
$ws.Range("B192").Value = "237675555508"
$ws.Range("B192").Style = "Normal"
$ws.Range("C192").Value = "LA NEGRESSE LTDLA CBOX R0 OKALA NTCHAGOU Ernest"
$ws.Range("D192").Value = 697
$ws.Range("A193").NumberFormat = "@"
$ws.Range("A193").Value = "2026-02-15 10:45:42"
$ws.Range("A193").Style = "Normal"
$ws.Range("B193").NumberFormat = "@"
$ws.Range("B193").Value = "237675557616"
$ws.Range("B193").Style = "Normal"
$ws.Range("C193").Value = "LUCIENNE FOTSING TSINGOUM"
$ws.Range("D193").Value = 568707
$ws.Range("A194").NumberFormat = "@"
$ws.Range("A194").Value = "2026-02-15 01:37:05"
$ws.Range("A194").Style = "Normal"
$ws.Range("B194").NumberFormat = "@"
$ws.Range("B194").Value = "237676301061"
$ws.Range("B194").Style = "Normal"
$ws.Range("C194").Value = "NGO TONYE ELISE AUDREY STYLE.COM"
$ws.Range("D194").Value = 75897
$ws.Range("A195").NumberFormat = "@"
$ws.Range("A195").Value = "2026-02-15 02:00:25"
$ws.Range("A195").Style = "Normal"
$ws.Range("B195").NumberFormat = "@"
$ws.Range("B195").Value = "237677939039"
$ws.Range("B195").Style = "Normal"
$ws.Range("C195").Value = "ROMARIC TRESOR TCHOUNKEU MBAKOP ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Range("D195").Value = 442
$ws.Range("A196").NumberFormat = "@"
$ws.Range("A196").Value = "2026-02-15 11:43:32"
$ws.Range("A196").Style = "Normal"
$ws.Range("B196").NumberFormat = "@"
$ws.Range("B196").Value = "237678225987"
$ws.Range("B196").Style = "Normal"
$ws.Range("C196").Value = "FINGO FOTO SOPHIE ETS LE CONTENT"
$ws.Range("D196").Value = 19174
$ws.Range("A197").NumberFormat = "@"
$ws.Range("A197").Value = "2026-02-15 14:34:12"
$ws.Range("A197").Style = "Normal"
$ws.Range("B197").NumberFormat = "@"
$ws.Range("B197").Value = "237679553674"
$ws.Range("B197").Style = "Normal"
$ws.Range("C197").Value = "DESIRE MAGINZANG MBOUEZKO"
$ws.Range("D197").Value = 202671
$ws.Range("A198").NumberFormat = "@"
$ws.Range("A198").Value = "2026-02-15 16:08:24"
$ws.Range("A198").Style = "Normal"
$ws.Range("B198").NumberFormat = "@"
$ws.Range("B198").Value = "237679654555"
$ws.Range("B198").Style = "Normal"
$ws.Range("C198").Value = "PAULINE NGUELEMO"
$ws.Range("D198").Value = 11435
$ws.Range("A199").NumberFormat = "@"
$ws.Range("A199").Value = "2026-02-13 02:55:19"
$ws.Range("A199").Style = "Normal"
$ws.Range("B199").NumberFormat = "@"
$ws.Range("B199").Value = "237679789713"
$ws.Range("B199").Style = "Normal"
$ws.Range("C199").Value = "ETS LE CONTENT 68"
$ws.Range("D199").Value = 48
$ws.Range("A200").NumberFormat = "@"
$ws.Range("A200").Value = "2026-02-15 14:54:22"
$ws.Range("A200").Style = "Normal"
$ws.Range("B200").NumberFormat = "@"
$ws.Range("B200").Value = "237681589841"
$ws.Range("B200").Style = "Normal"
$ws.Range("C200").Value = "LEINTENG ROSE MARY"
$ws.Range("D200").Value = 166469
$ws.Range("A201").NumberFormat = "@"
$ws.Range("A201").Value = "2026-02-15 11:31:27"
$ws.Range("A201").Style = "Normal"
$ws.Range("B201").NumberFormat = "@"
$ws.Range("B201").Value = "237681676445"
$ws.Range("B201").Style = "Normal"
$ws.Range("C201").Value = "MELANIE NGAFFO"
$ws.Range("D201").Value = 232346
$ws.Range("A202").NumberFormat = "@"
$ws.Range("A202").Value = "2026-02-14 11:23:00"
$ws.Range("A202").Style = "Normal"
$ws.Range("B202").NumberFormat = "@"
$ws.Range("B202").Value = "237681678622"
$ws.Range("B202").Style = "Normal"
$ws.Range("C202").Value = "FOKGO BRIGITTE ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Range("D202").Value = 90015
$ws.Range("A203").NumberFormat = "@"
$ws.Range("A203").Value = "2026-02-15 16:44:27"
$ws.Range("A203").Style = "Normal"
$ws.Range("B203").NumberFormat = "@"
$ws.Range("B203").Value = "237683815311"
$ws.Range("B203").Style = "Normal"
$ws.Range("C203").Value = "ETS TIN-GLOBALCOMM ZEBAZE TSEBAZE LAURA"
$ws.Range("D203").Value = 32404
$ws.Range("A204").NumberFormat = "@"
$ws.Range("A204").Value = "2026-02-15 11:26:37"
$ws.Range("A204").Style = "Normal"
$ws.Range("B204").NumberFormat = "@"
$ws.Range("B204").Value = "237651533411"
$ws.Range("B204").Style = "Normal"
$ws.Range("C204").Value = "CLAUDE LANDRY DJOUMSI NZUSSING"
$ws.Range("D204").Value = 35554
$ws.Range("A205").NumberFormat = "@"
$ws.Range("A205").Value = "2026-02-15 13:52:17"
$ws.Range("A205").Style = "Normal"
$ws.Range("B205").NumberFormat = "@"
$ws.Range("B205").Value = "237651843112"
$ws.Range("B205").Style = "Normal"
$ws.Range("C205").Value = "NLOGA NGO SIPORAH FELICITE ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Range("D205").Value = 5571
$ws.Range("A206").NumberFormat = "@"
$ws.Range("A206").Value = "2026-02-15 12:37:52"
$ws.Range("A206").Style = "Normal"
$ws.Range("B206").NumberFormat = "@"
$ws.Range("B206").Value = "237652297747"
$ws.Range("B206").Style = "Normal"
$ws.Range("C206").Value = "GOUAMPOUM MIREILLE FLORE ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Range("D206").Value = 159435
$ws.Range("A207").NumberFormat = "@"
$ws.Range("A207").Value = "2026-02-15 16:06:07"
$ws.Range("A207").Style = "Normal"
$ws.Range("B207").NumberFormat = "@"
$ws.Range("B207").Value = "237653316656"
$ws.Range("B207").Style = "Normal"
$ws.Range("C207").Value = "deschance zambou dontsop"
$ws.Range("D207").Value = 454540
$ws.Range("A208").NumberFormat = "@"
$ws.Range("A208").Value = "2026-02-15 12:46:28"
$ws.Range("A208").Style = "Normal"
$ws.Range("B208").NumberFormat = "@"
$ws.Range("B208").Value = "237671290825"
$ws.Range("B208").Style = "Normal"
$ws.Range("C208").Value = "LONGMENE FLORIDE NINA ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Range("D208").Value = 83290
$ws.Range("A209").NumberFormat = "@"
$ws.Range("A209").Value = "2026-02-14 19:34:23"
$ws.Range("A209").Style = "Normal"
$ws.Range("B209").NumberFormat = "@"
$ws.Range("B209").Value = "237672276931"
$ws.Range("B209").Style = "Normal"
$ws.Range("C209").Value = "MIREILLE MAGOLACK FOMEKONG"
$ws.Range("D209").Value = 14
$ws.Range("A210").NumberFormat = "@"
$ws.Range("A210").Value = "2026-02-15 09:40:30"
$ws.Range("A210").Style = "Normal"
$ws.Range("B210").NumberFormat = "@"
$ws.Range("B210").Value = "237674580187"
$ws.Range("B210").Style = "Normal"
$ws.Range("C210").Value = "ROSETTE SIGHOM"
$ws.Range("D210").Value = 21475
$ws.Range("A211").NumberFormat = "@"
$ws.Range("A211").Value = "2026-02-15 16:01:43"
$ws.Range("A211").Style = "Normal"
$ws.Range("B211").NumberFormat = "@"
$ws.Range("B211").Value = "237674673359"
$ws.Range("B211").Style = "Normal"
$ws.Range("C211").Value = "NDE DONATUS ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Range("D211").Value = 857779
$ws.Range("A212").NumberFormat = "@"
$ws.Range("A212").Value = "2026-02-15 18:01:11"
$ws.Range("A212").Style = "Normal"
$ws.Range("B212").NumberFormat = "@"
$ws.Range("B212").Value = "237674769580"
$ws.Range("B212").Style = "Normal"
$ws.Range("C212").Value = "VAKSA BELLO"
$ws.Range("D212").Value = 19533
$ws.Range("A213").NumberFormat = "@"
$ws.Range("A213").Value = "2026-02-15 18:24:11"
$ws.Range("A213").Style = "Normal"
$ws.Range("B213").NumberFormat = "@"
$ws.Range("B213").Value = "237674926472"
$ws.Range("B213").Style = "Normal"
$ws.Range("C213").Value = "DYLAN LECANAL MOLUH"
$ws.Range("D213").Value = 102388
$ws.Range("A214").NumberFormat = "@"
$ws.Range("A214").Value = "2026-02-15 03:43:12"
$ws.Range("A214").Style = "Normal"
$ws.Range("B214").NumberFormat = "@"
$ws.Range("B214").Value = "237675831509"
$ws.Range("B214").Style = "Normal"
$ws.Range("C214").Value = "NGAMOUN NAFISSATOU ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Range("D214").Value = 7307
$ws.Range("A215").NumberFormat = "@"
$ws.Range("A215").Value = "2026-02-15 15:05:06"
$ws.Range("A215").Style = "Normal"
$ws.Range("B215").NumberFormat = "@"
$ws.Range("B215").Value = "237680857383"
$ws.Range("B215").Style = "Normal"
$ws.Range("C215").Value = "MOUNIKIEL TECLAIRE ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Range("D215").Value = 570230
$ws.Range("A216").NumberFormat = "@"
$ws.Range("A216").Value = "2026-02-15 10:45:26"
$ws.Range("A216").Style = "Normal"
$ws.Range("B216").NumberFormat = "@"
$ws.Range("B216").Value = "237681180496"
$ws.Range("B216").Style = "Normal"
$ws.Range("C216").Value = "ALEX NGOUO YOUNDA"
$ws.Range("D216").Value = 4972
$ws.Range("A217").NumberFormat = "@"
$ws.Range("A217").Value = "2026-02-15 14:09:39"
$ws.Range("A217").Style = "Normal"
$ws.Range("B217").NumberFormat = "@"
$ws.Range("B217").Value = "237681299829"
$ws.Range("B217").Style = "Normal"
$ws.Range("C217").Value = "NDEBI MEDARD DESIRE ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Range("D217").Value = 28395
$ws.Range("A218").NumberFormat = "@"
$ws.Range("A218").Value = "2026-02-15 16:10:11"
$ws.Range("A218").Style = "Normal"
$ws.Range("B218").NumberFormat = "@"
$ws.Range("B218").Value = "237681657939"
$ws.Range("B218").Style = "Normal"
$ws.Range("C218").Value = "ETS MOBILE FINANCIAL SERVICES MFS LTDLA CBOX R1 MOHA CHAIBOU"
$ws.Range("D218").Value = 100903
$ws.Range("A219").NumberFormat = "@"
$ws.Range("A219").Value = "2026-02-15 14:33:09"
$ws.Range("A219").Style = "Normal"
$ws.Range("B219").NumberFormat = "@"
$ws.Range("B219").Value = "237675944533"
$ws.Range("B219").Style = "Normal"
$ws.Range("C219").Value = "LONGA ELDRINE ELYSEE ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Range("D219").Value = 189399
